# Hotfix: Fri Nov 15 16:48:15 RTZ 2024
# Adds a "last modified" timestamp column (D) and a (currently empty)
# trailing column (E) to every data row of the "Bash" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bash")

$lastRow = 81
$timestamp = "2024-11-15 13:18:16"

# Stamp every data row (1-81) in column D with the same "last run" timestamp.
$ws.Range("D1:D" + $lastRow).Value = $timestamp

# Materialize column E as existing-but-empty cells for every row so the
# sheet's used range (dimension) extends through E81, matching the
# placeholder column added alongside the timestamps. A plain Value=""
# assignment does not create a cell at all, so we touch a formatting
# property instead, which forces the (blank) cell to be written out.
$ws.Range("E1:E" + $lastRow).WrapText = $false
